$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Wed Feb 15 12:44:39 EST 2023"
$ws.Range("B3").Value = "Wed Feb 15 12:44:48 EST 2023"
$ws.Range("B4").Value = "Wed Feb 15 12:44:59 EST 2023"
